$wb = $excel.ActiveWorkbook
$timestamp = "2025-11-07 02:49:21"

$ws = $wb.Worksheets.Item(2)

# Update data values for rows with changed stats
# Row 7
$ws.Range("C7").Value = 14
$ws.Range("D7").Value = 187
$ws.Range("E7").Value = 80
$ws.Range("F7").Value = 107
$ws.Range("G7").Value = 13.36
$ws.Range("H7").Value = 5.71
$ws.Range("I7").Value = 7.64
$ws.Range("J7").Value = 40
$ws.Range("K7").Value = 36
# Row 13
$ws.Range("C13").Value = 9
$ws.Range("D13").Value = 145
$ws.Range("E13").Value = 80
$ws.Range("F13").Value = 65
$ws.Range("G13").Value = 16.11
$ws.Range("H13").Value = 8.89
$ws.Range("I13").Value = 7.22
$ws.Range("J13").Value = 40
$ws.Range("K13").Value = 30
# Row 16
$ws.Range("C16").Value = 22
$ws.Range("D16").Value = 440
$ws.Range("E16").Value = 213
$ws.Range("F16").Value = 227
$ws.Range("G16").Value = 20
$ws.Range("H16").Value = 9.68
$ws.Range("I16").Value = 10.32
$ws.Range("J16").Value = 79
$ws.Range("K16").Value = 81
# Row 17
$ws.Range("C17").Value = 15
$ws.Range("D17").Value = 249
$ws.Range("E17").Value = 90
$ws.Range("F17").Value = 159
$ws.Range("G17").Value = 16.6
$ws.Range("H17").Value = 6
$ws.Range("I17").Value = 10.6
$ws.Range("J17").Value = 45
$ws.Range("K17").Value = 62
$ws.Range("M17").Value = 1
$ws.Range("Y17").Value = 2
# Row 18
$ws.Range("C18").Value = 22
$ws.Range("D18").Value = 330
$ws.Range("E18").Value = 154
$ws.Range("F18").Value = 176
$ws.Range("G18").Value = 15
$ws.Range("H18").Value = 7
$ws.Range("I18").Value = 8
$ws.Range("J18").Value = 67
$ws.Range("K18").Value = 83
$ws.Range("Y18").Value = 6
# Row 20
$ws.Range("C20").Value = 21
$ws.Range("D20").Value = 363
$ws.Range("E20").Value = 148
$ws.Range("F20").Value = 215
$ws.Range("G20").Value = 17.29
$ws.Range("H20").Value = 7.05
$ws.Range("I20").Value = 10.24
$ws.Range("J20").Value = 69
$ws.Range("K20").Value = 80
$ws.Range("Y20").Value = 4
# Row 21
$ws.Range("C21").Value = 18
$ws.Range("D21").Value = 269
$ws.Range("E21").Value = 118
$ws.Range("F21").Value = 151
$ws.Range("G21").Value = 14.94
$ws.Range("H21").Value = 6.56
$ws.Range("I21").Value = 8.39
$ws.Range("J21").Value = 49
$ws.Range("K21").Value = 63
# Row 23
$ws.Range("C23").Value = 13
$ws.Range("D23").Value = 175
$ws.Range("E23").Value = 62
$ws.Range("F23").Value = 113
$ws.Range("G23").Value = 13.46
$ws.Range("H23").Value = 4.77
$ws.Range("I23").Value = 8.69
$ws.Range("J23").Value = 31
$ws.Range("K23").Value = 44
$ws.Range("M23").Value = 1
$ws.Range("Y23").Value = 4

# Update as_of_utc timestamp for every data row (2-26)
for ($r = 2; $r -le 26; $r++) {
    $ws.Range("AA" + $r).Value = $timestamp
}

$ws = $wb.Worksheets.Item(3)

# Update data values for rows with changed stats
# Row 3
$ws.Range("C3").Value = 20
$ws.Range("D3").Value = 298
$ws.Range("E3").Value = 142
$ws.Range("F3").Value = 156
$ws.Range("G3").Value = 14.9
$ws.Range("H3").Value = 7.1
$ws.Range("I3").Value = 7.8
$ws.Range("J3").Value = 71
$ws.Range("K3").Value = 63
# Row 12
$ws.Range("C12").Value = 19
$ws.Range("D12").Value = 336
$ws.Range("E12").Value = 159
$ws.Range("F12").Value = 177
$ws.Range("G12").Value = 17.68
$ws.Range("H12").Value = 8.37
$ws.Range("I12").Value = 9.32
$ws.Range("J12").Value = 72
$ws.Range("K12").Value = 81
$ws.Range("V12").Value = 10
# Row 14
$ws.Range("C14").Value = 22
$ws.Range("D14").Value = 378
$ws.Range("E14").Value = 188
$ws.Range("F14").Value = 190
$ws.Range("G14").Value = 17.18
$ws.Range("H14").Value = 8.55
$ws.Range("I14").Value = 8.64
$ws.Range("J14").Value = 94
$ws.Range("K14").Value = 90
$ws.Range("Y14").Value = 2
# Row 19
$ws.Range("C19").Value = 20
$ws.Range("D19").Value = 359
$ws.Range("E19").Value = 168
$ws.Range("F19").Value = 191
$ws.Range("G19").Value = 17.95
$ws.Range("H19").Value = 8.4
$ws.Range("I19").Value = 9.55
$ws.Range("J19").Value = 79
$ws.Range("K19").Value = 83
# Row 20
$ws.Range("C20").Value = 15
$ws.Range("D20").Value = 260
$ws.Range("E20").Value = 131
$ws.Range("F20").Value = 129
$ws.Range("G20").Value = 17.33
$ws.Range("H20").Value = 8.73
$ws.Range("I20").Value = 8.6
$ws.Range("J20").Value = 63
$ws.Range("K20").Value = 62
$ws.Range("Y20").Value = 4
# Row 21
$ws.Range("C21").Value = 24
$ws.Range("D21").Value = 511
$ws.Range("E21").Value = 213
$ws.Range("F21").Value = 298
$ws.Range("G21").Value = 21.29
$ws.Range("H21").Value = 8.88
$ws.Range("I21").Value = 12.42
$ws.Range("J21").Value = 99
$ws.Range("K21").Value = 119
$ws.Range("M21").Value = 2
$ws.Range("Y21").Value = 6
# Row 25
$ws.Range("C25").Value = 7
$ws.Range("D25").Value = 131
$ws.Range("E25").Value = 76
$ws.Range("F25").Value = 55
$ws.Range("G25").Value = 18.71
$ws.Range("H25").Value = 10.86
$ws.Range("I25").Value = 7.86
$ws.Range("J25").Value = 38
$ws.Range("K25").Value = 25

# Update as_of_utc timestamp for every data row (2-26)
for ($r = 2; $r -le 26; $r++) {
    $ws.Range("AA" + $r).Value = $timestamp
}
